$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity (C28) from 1 to 10
$ws.Range("C28").Value = 10

# Set unit price (G28) which was previously blank
$ws.Range("G28").Value = 0.02

# Update price total (H28) from 0.03 to 0.2
$ws.Range("H28").Value = 0.2

# Update sheet view: scroll position and active selection to I28
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I28").Select()
